# ----------------------------------------------------------------------
# reg_exceptional_holiday.xlsx -- "sierra leone master data"
#
# Re-purpose the Madagascar "Emergency Holiday" template: language
# becomes "eng", the holiday date is stored as plain text, three more
# registration centres (10002-10004) are added, and the whole data
# block is re-styled (bold Cambria header / italic Calibri body, new
# box borders, right/centre aligned columns).
# ----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- workbook level ------------------------------------------------
$excel.Iteration = $true
$excel.MaxChange = 0.0001

# ---- data ------------------------------------------------------------
# Row 2 & 3 already exist - update language, regcntr_id and the date
# (switched from a real date serial to the literal text "9/24/2023").
$ws.Range("C2:C5").NumberFormat = "@"

$ws.Range("A2").Value2 = "eng"
$ws.Range("B2").Value2 = 10001
$ws.Range("C2").Value2 = "9/24/2023"

$ws.Range("A3").Value2 = "eng"
$ws.Range("B3").Value2 = 10002
$ws.Range("C3").Value2 = "9/24/2023"

# Rows 4 & 5 are brand new registration centres with the same holiday.
$ws.Range("A4").Value2 = "eng"
$ws.Range("B4").Value2 = 10003
$ws.Range("C4").Value2 = "9/24/2023"
$ws.Range("D4").Value2 = "Emergency Holiday"
$ws.Range("E4").Value2 = "Emergency Holiday"
$ws.Range("F4").Value2 = $true

$ws.Range("A5").Value2 = "eng"
$ws.Range("B5").Value2 = 10004
$ws.Range("C5").Value2 = "9/24/2023"
$ws.Range("D5").Value2 = "Emergency Holiday"
$ws.Range("E5").Value2 = "Emergency Holiday"
$ws.Range("F5").Value2 = $true

$ws.Rows("2:5").RowHeight = 27
$ws.Rows("1:1").RowHeight = 15

Write-Host "data written"

# ---- fonts -------------------------------------------------------------
# Header font: bold Cambria 11pt (was plain Cambria 8pt).
$headerFont = $ws.Range("A1:F1").Font
$headerFont.Name = "Cambria"
$headerFont.Size = 11
$headerFont.Bold = $true
$headerFont.Italic = $false

# Body font: italic Calibri 10pt black -- now shared by every data
# column, including regcntr_id which used to be plain 9pt Arial.
$bodyFont = $ws.Range("A2:F5").Font
$bodyFont.Name = "Calibri"
$bodyFont.Size = 10
$bodyFont.Italic = $true
$bodyFont.Bold = $false
$bodyFont.Color = 0

Write-Host "fonts set"

# ---- alignment -----------------------------------------------------
$hdr = $ws.Range("A1:F1")
$hdr.HorizontalAlignment = -4108   # xlCenter
$hdr.VerticalAlignment = -4160     # xlTop
$hdr.WrapText = $false

$all = $ws.Range("A2:F5")
$all.WrapText = $true
$all.HorizontalAlignment = -4131   # xlGeneral
$all.VerticalAlignment = -4107     # xlBottom (default)

$rightAlign = $ws.Range("B2:C5")
$rightAlign.HorizontalAlignment = -4152   # xlRight

$centerAlign = $ws.Range("F2:F5")
$centerAlign.HorizontalAlignment = -4108  # xlCenter

Write-Host "alignment set"

# ---- fills ------------------------------------------------------------
# White solid fill highlighting regcntr_id and is_active columns.
$ws.Range("B2:B5").Interior.Color = 16777215
$ws.Range("F2:F5").Interior.Color = 16777215

Write-Host "fills set"

# ---- number formats ---------------------------------------------------
# Header cells over the "style" columns (B & F) apply the Text format,
# matching the column-level style used for those two columns.
$ws.Range("B1").NumberFormat = "@"
$ws.Range("F1").NumberFormat = "@"

Write-Host "numfmt set"

# ---- borders -----------------------------------------------------------
# Thin automatic-colour box around the header row.
$hdrB = $ws.Range("A1:F1").Borders
$hdrB.LineStyle = 1
$hdrB.Weight = 2
$hdrB.ColorIndex = -4105

# Row 2 (top data row): full medium box, left edge of B:F is a lighter
# grey because it touches the previous (A) cell, matching the original
# template's "interior seam" colouring.
$rowTop = $ws.Range("A2:F2")
$rowTop.Borders.LineStyle = 1
$rowTop.Borders.Weight = -4138
$rowTop.Borders.Color = 0
$ws.Range("B2:F2").Borders.Item(7).Color = 13421772

# Rows 3-5 (continuation rows): same as above, but the top edge is also
# the lighter grey "interior seam" colour.
$rowRest = $ws.Range("A3:F5")
$rowRest.Borders.LineStyle = 1
$rowRest.Borders.Weight = -4138
$rowRest.Borders.Color = 0
$ws.Range("A3:A5").Borders.Item(8).Color = 13421772
$ws.Range("B3:F5").Borders.Item(7).Color = 13421772
$ws.Range("B3:F5").Borders.Item(8).Color = 13421772

Write-Host "borders set"

# ---- columns -----------------------------------------------------------
# Drop the custom widths that used to widen hol_name/hol_reason.
$ws.Columns("D:E").ColumnWidth = 8.43

Write-Host "columns set"

# ---- view / selection ---------------------------------------------------
$ws.Range("J2").Select()
$excel.ActiveWindow.Zoom = 100
$excel.ActiveWindow.TabRatio = 500

Write-Host "view set"

# ---- page setup ----------------------------------------------------------
$ps = $ws.PageSetup
$ps.PaperSize = 9          # xlPaperA4
$ps.FirstPageNumber = 0
$ps.Orientation = 1        # xlPortrait
$ps.LeftMargin = $excel.InchesToPoints(0.75)
$ps.RightMargin = $excel.InchesToPoints(0.75)
$ps.TopMargin = $excel.InchesToPoints(1)
$ps.BottomMargin = $excel.InchesToPoints(1)
$ps.HeaderMargin = $excel.InchesToPoints(0.511805555555555)
$ps.FooterMargin = $excel.InchesToPoints(0.511805555555555)

Write-Host "page setup done"
